# Updating projet_list_all with project test set (3 projects per instrument)
# Adds a new "object_annotation_category" value in column Y for each data
# row of the "Data" sheet (the header row already has a value in Y1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 25).Value = "object_annotation_category"  # column Y
}
